# Update "想去人数" (interest count) column F values across all sheets
# per the recorded diff (gh-pages output regeneration).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 503
$ws.Range("F5").Value = 2345
$ws.Range("F6").Value = 9
$ws.Range("F7").Value = 66
$ws.Range("F8").Value = 78
$ws.Range("F9").Value = 1679
$ws.Range("F10").Value = 1679
$ws.Range("F11").Value = 1389
$ws.Range("F12").Value = 75
$ws.Range("F13").Value = 1433
$ws.Range("F16").Value = 789
$ws.Range("F17").Value = 1
$ws.Range("F19").Value = 131
$ws.Range("F20").Value = 7447
$ws.Range("F21").Value = 8355
$ws.Range("F26").Value = 96
$ws.Range("F34").Value = 1484
$ws.Range("F38").Value = 301
$ws.Range("F39").Value = 29
$ws.Range("F40").Value = 774
$ws.Range("F43").Value = 365
$ws.Range("F44").Value = 265
$ws.Range("F45").Value = 213
$ws.Range("F47").Value = 203
$ws.Range("F48").Value = 186
$ws.Range("F49").Value = 25

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 69
$ws.Range("F18").Value = 308

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 151

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 151
$ws.Range("F7").Value = 503
$ws.Range("F8").Value = 2345
$ws.Range("F9").Value = 9
$ws.Range("F10").Value = 66
$ws.Range("F11").Value = 78
$ws.Range("F12").Value = 1679
$ws.Range("F13").Value = 1679
$ws.Range("F15").Value = 75
$ws.Range("F16").Value = 1433
$ws.Range("F19").Value = 789
$ws.Range("F22").Value = 131
$ws.Range("F23").Value = 7447
$ws.Range("F24").Value = 7447
$ws.Range("F25").Value = 8356
$ws.Range("F28").Value = 96
$ws.Range("F37").Value = 301
$ws.Range("F38").Value = 29
$ws.Range("F41").Value = 774
$ws.Range("F45").Value = 365
$ws.Range("F46").Value = 265
$ws.Range("F47").Value = 213
$ws.Range("F48").Value = 203
$ws.Range("F49").Value = 186
$ws.Range("F50").Value = 308
$ws.Range("F51").Value = 25

